# Weekly update: insert a new daily price record at the top of the
# "Ajo" (garlic) data block for "Terminal La Palmera de La Serena",
# pushing the existing rows (221-275) down by one (to 222-276) and
# filling the freshly inserted row 221 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 221:275 down to 222:276, growing the used range to A1:R276.
$ws.Rows.Item(221).Insert()

# Populate the new row 221 with the new weekly record. Most fields repeat
# the values of the record that is now in row 222 (previously row 221);
# only the date (D) and volume (J) differ for the new entry.
$ws.Range("A221").Value = 8
$ws.Range("B221").Value = "Terminal La Palmera de La Serena"
$ws.Range("C221").Value = "Coquimbo"
$ws.Range("D221").Value = 44711
$ws.Range("E221").Value = 4
$ws.Range("F221").Value = 100112003
$ws.Range("G221").Value = "Ajo"
$ws.Range("H221").Value = "Chino"
$ws.Range("I221").Value = "Primera"
$ws.Range("J221").Value = 480
$ws.Range("K221").Value = 18000
$ws.Range("L221").Value = 19000
$ws.Range("M221").Value = 18500
$ws.Range("N221").Value = "`$/caja 10 kilos"
$ws.Range("O221").Value = "China"
$ws.Range("P221").Value = 1850
$ws.Range("Q221").Value = 10
$ws.Range("R221").Value = "Hortaliza"
